# DhalisMenu_cat.xlsx — "Add files via upload"
# Adds two new menu rows (Cream bell kesar badam milk / DEW20) to Sheet1,
# resizes the data columns to fit their content, and leaves the selection
# on the last cell entered, matching what Excel records after a manual
# data-entry + column AutoFit session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 54 ---------------------------------------------------------
$ws.Range("A54").Value = "Drinks"
$ws.Range("B54").Value = "Cream bell kesar badam milk"
$ws.Range("C54").Value = 0
$ws.Range("D54").Value = 20
$ws.Range("E54").Value = "Cream bell kesar badam milk.JPG"
$ws.Range("F54").Value = "Fast Food"

# --- New row 55 ---------------------------------------------------------
$ws.Range("A55").Value = "Drinks"
$ws.Range("B55").Value = "DEW20"
$ws.Range("C55").Value = 0
$ws.Range("D55").Value = 20
$ws.Range("E55").Value = "DEW20.JPG"
$ws.Range("F55").Value = "Fast Food"

# --- Resize the columns that hold the menu text/image data so the new,
#     longer entries are fully visible (mirrors double-clicking each
#     column border to best-fit the content). ---------------------------
$ws.Range("A:A").ColumnWidth = 7.983072916666667
$ws.Range("B:B").ColumnWidth = 65.61979166666667
$ws.Range("C:C").ColumnWidth = 3.3463541666666665
$ws.Range("E:E").ColumnWidth = 52.166666666666664
$ws.Range("F:F").ColumnWidth = 7.983072916666667

# --- Leave the selection where data entry finished ----------------------
$ws.Range("F55").Select()
